$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "29.365.69"
$ws.Range("E2").Value = "  -0.16%  "
$ws.Range("D3").Value = "1.845.11"
$ws.Range("E3").Value = "  -0.32%  "
$ws.Range("E4").Value = "  +0.01%  "
$ws.Range("D5").Value = "'239.00"
$ws.Range("E5").Value = "  -1.34%  "
$ws.Range("D6").Value = "'0.6292"
$ws.Range("E6").Value = "  +0.00%  "
$ws.Range("E7").Value = "  +0.01%  "
$ws.Range("D8").Value = "'0.07551"
$ws.Range("E8").Value = "  -1.58%  "
$ws.Range("D9").Value = "'0.2944"
$ws.Range("E9").Value = "  -1.08%  "
$ws.Range("D10").Value = "'24.56"
$ws.Range("E10").Value = "  +0.37%  "
$ws.Range("D11").Value = "'0.07693"
$ws.Range("E11").Value = "  -0.36%  "
$ws.Range("D12").Value = "1.838.93"
$ws.Range("E12").Value = "  -6.18%  "
$ws.Range("E13").Value = "  -0.55%  "
$ws.Range("D14").Value = "'0.6782"
$ws.Range("E14").Value = "  -1.60%  "
$ws.Range("D15").Value = "'0.00001023"
$ws.Range("E15").Value = "  +2.47%  "
$ws.Range("E16").Value = "  -0.10%  "
$ws.Range("D17").Value = "2.090.34"
$ws.Range("E17").Value = "  -4.69%  "
$ws.Range("D18").Value = "'6.129"
$ws.Range("E18").Value = "  -0.96%  "
$ws.Range("D19").Value = "29.405.41"
$ws.Range("E19").Value = "  -0.50%  "
$ws.Range("D20").Value = "'228.18"
$ws.Range("E20").Value = "  -2.06%  "
$ws.Range("E21").Value = "  -1.21%  "
$ws.Range("D22").Value = "'1.000"
$ws.Range("E22").Value = "  +0.05%  "
$ws.Range("E23").Value = "  -2.75%  "
$ws.Range("E24").Value = "  +0.07%  "
$ws.Range("E25").Value = "  +1.26%  "
$ws.Range("E26").Value = "  -0.04%  "
$ws.Range("D27").Value = "'8.349"
$ws.Range("E27").Value = "  -1.64%  "
$ws.Range("D28").Value = "'17.61"
$ws.Range("E28").Value = "  -0.31%  "
$ws.Range("D29").Value = "'1.459"
$ws.Range("E29").Value = "  -1.31%  "
$ws.Range("D30").Value = "'1.264"
$ws.Range("E30").Value = "  +0.58%  "
$ws.Range("D31").Value = "'0.05632"
$ws.Range("E31").Value = "  -2.68%  "
$ws.Range("E32").Value = "  -0.24%  "
$ws.Range("D33").Value = "'4.034"
$ws.Range("E33").Value = "  +0.27%  "
$ws.Range("D34").Value = "'1.834"
$ws.Range("E34").Value = "  -2.62%  "
$ws.Range("E35").Value = "  -0.57%  "
$ws.Range("D36").Value = "'0.7113"
$ws.Range("E36").Value = "  -1.68%  "
$ws.Range("D37").Value = "'2.592"
$ws.Range("E37").Value = "  +0.21%  "
$ws.Range("D38").Value = "1.239.76"
$ws.Range("E38").Value = "  -0.84%  "
$ws.Range("D39").Value = "'0.01807"
$ws.Range("E39").Value = "  -0.11%  "
$ws.Range("D40").Value = "'2.768"
$ws.Range("E40").Value = "  -1.01%  "
$ws.Range("D41").Value = "'6.228"
$ws.Range("E41").Value = "  +2.23%  "
$ws.Range("D42").Value = "'0.9010"
$ws.Range("E42").Value = "  -0.90%  "
$ws.Range("D44").Value = "'101.85"
$ws.Range("E44").Value = "  +0.31%  "
$ws.Range("D45").Value = "'65.86"
$ws.Range("E45").Value = "  -2.99%  "
$ws.Range("D46").Value = "'7.101"
$ws.Range("E46").Value = "  -2.57%  "
$ws.Range("D47").Value = "'0.00000000117"
$ws.Range("E47").Value = "  -2.65%  "
$ws.Range("D48").Value = "'0.3993"
$ws.Range("E48").Value = "  -0.98%  "
$ws.Range("D49").Value = "'1.683"
$ws.Range("E49").Value = "  -0.88%  "
$ws.Range("E50").Value = "  -2.77%  "
$ws.Range("D51").Value = "'0.1117"
$ws.Range("E51").Value = "  -0.43%  "
